$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FA bugs")

# The "Highlight Sites..." block (rows 5:8) needs to move up two rows to
# become rows 3:6, because the old row-3 item ("Finish setting download
# headers") has been resolved and relocated further down the sheet.
# Deleting the old row 3 plus its blank spacer row (row 4) shifts
# everything below up by exactly two rows, which reproduces that move
# (and all the other rows below it) in one shot.
$ws.Rows("3:4").Delete()

# That same delete shifted the old row 26 ("selection box issues") block
# up to row 24, which needs to stay at row 26. Insert two blank rows
# back in before it to restore rows 26+ to their original positions.
$ws.Rows("24:25").Insert(-4121)

# Row 25 must end up completely empty (no leftover formatting), while
# row 24 receives the relocated "Finish setting download headers" item.
$ws.Range("A25").Clear()

# Give A24:C24 the same (grey, wrap-text) look already used by the
# "Highlight Sites" row (now at A3) by copying its format, then fill in
# the relocated text.
$ws.Range("A3").Copy()
$ws.Range("A24:C24").PasteSpecial(-4122)
$ws.Range("A24").Value = "Finish setting download headers"
$ws.Range("B24").Value = "y"
$ws.Range("C24").Value = "y"

# Move the selection cursor to A5, matching the saved view state.
$ws.Activate()
$ws.Range("A5").Select()
